$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ13625108",
    "summ13705841",
    "summ13803991",
    "summ13887328",
    "summ13973609",
    "summ14057919",
    "summ14141135",
    "summ14223715",
    "summ14306943",
    "summ14392891",
    "summ14475729",
    "summ14556999",
    "summ14656228",
    "summ14740442",
    "summ14823341",
    "summ14906170",
    "summ14990318",
    "summ15088180",
    "summ15185528",
    "summ15271582",
    "summ15361187",
    "summ15441545",
    "summ15523053",
    "summ15607062",
    "summ15691542",
    "summ15876058",
    "summ15956808",
    "summ16038129",
    "summ16125988",
    "summ16206919",
    "summ16290418",
    "summ16372176",
    "summ16457066",
    "summ16541417",
    "summ16623057",
    "summ16707117",
    "summ16804822",
    "summ16899041",
    "summ16975761",
    "summ17057102",
    "summ17159231",
    "summ17242164",
    "summ17353246",
    "summ17439553",
    "summ17522428",
    "summ17606254",
    "summ17689368",
    "summ17774029",
    "summ17856535",
    "summ17940438"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

